$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44253
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 1200
$ws.Range("K4").Value = 270
$ws.Range("L4").Value = 280
$ws.Range("M4").Value = 275
$ws.Range("P4").Value = 275

# Row 5
$ws.Range("D5").Value = 44229
$ws.Range("I5").Value = "Primera"
$ws.Range("K5").Value = 230
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = 240
$ws.Range("P5").Value = 240

# Row 6
$ws.Range("D6").Value = 44547
$ws.Range("K6").Value = 350
$ws.Range("L6").Value = 370
$ws.Range("M6").Value = 360
$ws.Range("P6").Value = 360
